# Insert a new weekly data row at row 18 (shifting the existing rows 18-167
# down to 19-168) and populate it with the new record, matching the
# metadata of the row that previously occupied position 18 (which is now
# row 19) but with the new date and volume for the new week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 18..167 down by one row.
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Range("A18").Value = 4
$ws.Range("B18").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C18").Value = "Los Lagos"
$ws.Range("D18").Value = 44532
$ws.Range("E18").Value = 10
$ws.Range("F18").Value = 100112017
$ws.Range("G18").Value = "Apio"
$ws.Range("H18").Value = "Americana (o)"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = 12000
$ws.Range("N18").Value = "$/docena de matas"
$ws.Range("O18").Value = "Región de Coquimbo"
$ws.Range("P18").Value = 2000
$ws.Range("Q18").Value = 6
$ws.Range("R18").Value = "Hortaliza"
